# Apply the StructureDefinition-episode-allowed-amount-rx.xlsx update:
#  - bump the IG version metadata (Metadata sheet)
#  - replace the stray duplicate "Contact" row with a proper Jurisdiction row
#  - give the top-level Extension row a real Short/Definition (Elements sheet)

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second, useless "Contact" / "No display for
# ContactDetail" pair; turn it into the Jurisdiction row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was an exact duplicate of the old row 10 ("Contact" / "No
# display for ContactDetail") - remove it outright, shifting
# Description/Purpose/.../Context up by one row.
$meta.Rows.Item(11).Delete()

# ---- Elements sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Top-level Extension row: give it the real Short/Definition text
# instead of the generic placeholder strings.
$elements.Range("K2").Value = "Episode Allowed Amount Rx (USD)"
$elements.Range("L2").Value = "Allowed amount per episode for the medications, in USD"
